$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "asdsad"
$ws.Range("C1").Value = "asd"

$ws.Range("C1").Select()
